$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$savedStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.996.40"
$ws.Range("D2").Style = $savedStyle
$ws.Range("E2").Value = "  +2.39%  "
$savedStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.999.69"
$ws.Range("D3").Style = $savedStyle
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  +0.16%  "
$savedStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.60"
$ws.Range("D5").Style = $savedStyle
$ws.Range("E5").Value = "  +1.32%  "
$savedStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.04"
$ws.Range("D6").Style = $savedStyle
$ws.Range("E6").Value = "  +3.41%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +1.18%  "
$savedStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.990.28"
$ws.Range("D9").Style = $savedStyle
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("E10").Value = "  +3.46%  "
$savedStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.23"
$ws.Range("D11").Style = $savedStyle
$ws.Range("E11").Value = "  +8.18%  "
$ws.Range("E12").Value = "  +1.94%  "
$savedStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000230"
$ws.Range("D13").Style = $savedStyle
$ws.Range("E13").Value = "  +3.18%  "
$savedStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.58"
$ws.Range("D14").Style = $savedStyle
$ws.Range("E14").Value = "  +2.12%  "
$ws.Range("E15").Value = "  +2.06%  "
$savedStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.495.35"
$ws.Range("D16").Style = $savedStyle
$ws.Range("E16").Value = "  +1.91%  "
$savedStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.29"
$ws.Range("D17").Style = $savedStyle
$ws.Range("E17").Value = "  +7.28%  "
$savedStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.994.57"
$ws.Range("D18").Style = $savedStyle
$ws.Range("E18").Value = "  +1.67%  "
$savedStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "58.984.22"
$ws.Range("D19").Style = $savedStyle
$ws.Range("E19").Value = "  +2.56%  "
$savedStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "428.22"
$ws.Range("D20").Style = $savedStyle
$ws.Range("E20").Value = "  +2.62%  "
$savedStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.72"
$ws.Range("D21").Style = $savedStyle
$ws.Range("E21").Value = "  +4.45%  "
$savedStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.722"
$ws.Range("D22").Style = $savedStyle
$ws.Range("E22").Value = "  +5.66%  "
$ws.Range("E23").Value = "  +2.00%  "
$savedStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.33"
$ws.Range("D24").Style = $savedStyle
$ws.Range("E24").Value = "  +2.45%  "
$savedStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.55"
$ws.Range("D25").Style = $savedStyle
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("E26").Value = "  +0.06%  "
$savedStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = $savedStyle
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("E28").Value = "  +10.45%  "
$ws.Range("E29").Value = "  +2.34%  "
$savedStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.79"
$ws.Range("D30").Style = $savedStyle
$ws.Range("E30").Value = "  +3.41%  "
$savedStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.75"
$ws.Range("D31").Style = $savedStyle
$ws.Range("E31").Value = "  +2.20%  "
$savedStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.04"
$ws.Range("D32").Style = $savedStyle
$ws.Range("E32").Value = "  -0.35%  "
$savedStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0993"
$ws.Range("D33").Style = $savedStyle
$ws.Range("E33").Value = "  -2.14%  "
$ws.Range("E34").Value = "  +6.40%  "
$ws.Range("E35").Value = "  +5.53%  "
$savedStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0757"
$ws.Range("D36").Style = $savedStyle
$ws.Range("E36").Value = "  +9.67%  "
$savedStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.10"
$ws.Range("D37").Style = $savedStyle
$ws.Range("E37").Value = "  -0.58%  "
$savedStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.71"
$ws.Range("D38").Style = $savedStyle
$ws.Range("E38").Value = "  +0.11%  "
$savedStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.67"
$ws.Range("D39").Style = $savedStyle
$ws.Range("E39").Value = "  +2.47%  "
$savedStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.73"
$ws.Range("D40").Style = $savedStyle
$ws.Range("E40").Value = "  +7.05%  "
$savedStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "398.57"
$ws.Range("D41").Style = $savedStyle
$ws.Range("E41").Value = "  +5.49%  "
$ws.Range("E42").Value = "  +0.62%  "
$savedStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.753.69"
$ws.Range("D43").Style = $savedStyle
$ws.Range("E43").Value = "  +3.30%  "
$ws.Range("E44").Value = "  -0.66%  "
$ws.Range("E45").Value = "  +4.79%  "
$savedStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "35.43"
$ws.Range("D46").Style = $savedStyle
$ws.Range("E46").Value = "  +25.17%  "
$savedStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.42"
$ws.Range("D48").Style = $savedStyle
$ws.Range("E48").Value = "  +0.69%  "
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("E50").Value = "  +0.57%  "
$savedStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.39"
$ws.Range("D51").Style = $savedStyle
$ws.Range("E51").Value = "  -0.20%  "
